# Apply the "everyday" management workbook edit:
# - Add a new shared string / cell value in sheet1 at A11
# - Move the active selection from A10 to B10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("everyday")

# Add the new note in row 11 (row 10 stays empty, matching the diff which
# keeps dimension A1:A11 but has no row 10 entry in sheetData).
$ws.Range("A11").Value = "工作时间也会成为绩效考核的一部分，虽然比例低，但是会考虑。"

# Update the active selection to B10, as reflected in the diff.
$ws.Range("B10").Select()
